# Weekly update: insert one new data row for "Perejil" at Vega Central
# Mapocho de Santiago, shifting all existing rows from 338 downward by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 338 (Excel shifts rows 338:421 -> 339:422,
# and the sheet's used range / dimension grows to A1:R422 automatically).
$ws.Rows.Item(338).Insert()

# Populate the newly inserted row 338 with the new weekly observation.
$ws.Range("A338").Value = 9
$ws.Range("B338").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C338").Value = "Metropolitana"
$ws.Range("D338").Value2 = 44785
$ws.Range("E338").Value = 13
$ws.Range("F338").Value = 100112044
$ws.Range("G338").Value = "Perejil"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 70
$ws.Range("K338").Value = 13000
$ws.Range("L338").Value = 15000
$ws.Range("M338").Value = 14000
$ws.Range("N338").Value = "$/docena de atados"
$ws.Range("O338").Value = "Región Metropolitana"
$ws.Range("P338").Value = 4667
$ws.Range("Q338").Value = 3
$ws.Range("R338").Value = "Hortaliza"
